$wb = $excel.ActiveWorkbook

# --- Sheet1: move the selection (cosmetic navigation left by the author) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("B41").Select() | Out-Null

# --- Updated-Todo: no content change, just loses the "last active" tab flag
#     automatically once a new sheet becomes active below ---
$ws2 = $wb.Worksheets.Item("Updated-Todo")

# --- Add the new "Update to-do 08-19-2018" sheet at the very end ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Update to-do 08-19-2018"

# Type the new to-do items in the same order the author actually typed them
# (this keeps the shared-string table ordering identical to the source file):
# first the tail of the list (rows 7-14), then the newly-inserted head
# (rows 1-5, typed in a slightly out-of-order fashion), then finally row 6,
# which re-uses a todo item already present elsewhere in the workbook.
$ws3.Cells.Item(7, 1).Value = "Think about new title on Plan a Trip Page"
$ws3.Cells.Item(8, 1).Value = "Think about splitting up locations and activities"
$ws3.Cells.Item(9, 1).Value = "Pull times from locations in order to correctly generate itinerary based on when they're open"
$ws3.Cells.Item(10, 1).Value = "Dynamic clicking on itinerary"
$ws3.Cells.Item(11, 1).Value = "Web spiders to parse information for top search locations to put into database for faster loading"
$ws3.Cells.Item(12, 1).Value = "Data analytics for pulling information"
$ws3.Cells.Item(13, 1).Value = "Database of items necessary based on activities in your trip"
$ws3.Cells.Item(14, 1).Value = "Generate itinerary creates a new page where the plan is full page with a map on the bottom with all of the directions from place to place where you can click on legs on either the map or the itinerary to highlight them, see the travel time, change mode of transportation, etc."

$ws3.Cells.Item(2, 1).Value = "Create database that is capable  of saving user names / passwords"
$ws3.Cells.Item(5, 1).Value = "Create database for itineraries"
$ws3.Cells.Item(4, 1).Value = "Back end to calculate the itinerary"
$ws3.Cells.Item(3, 1).Value = "Clickable google maps directions based on the calculated itinerary page"
$ws3.Cells.Item(1, 1).Value = "Create Back End (Stores Username / Password, Calculates Itinerary, Stores Itinerary, Shows Itinerary, Parses Data & Grabs data from internet etc.)"

$ws3.Cells.Item(6, 1).Value = "Be able to search along routes after itinerary is made for possible crags / gas stations"

# Re-apply the same (pre-existing) cell formatting the source workbook uses
# for these two particular rows (row 6 reuses the "search along routes"
# item's look, row 14 reuses the final highlighted-row look from Updated-Todo).
$ws1.Range("A22").Copy() | Out-Null
$ws3.Range("A6").PasteSpecial(-4122) | Out-Null

$ws2.Range("B27").Copy() | Out-Null
$ws3.Range("A14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Park the selection where the author left it on the new sheet.
$ws3.Range("L5").Select() | Out-Null
